# canteen:report generation timeout increased
# Refresh "This Month Report" sheet1 with the latest order data pulled
# from the canteen DB: rows are re-sorted to match the current query
# result order and three newly-placed orders are appended at the end.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

function Set-RowData {
    param($r, $orderId, $customerName, $totalPrice, $quantity, $paymentMode)
    $ws.Range("A$r").Value = $orderId
    $ws.Range("B$r").Value = $customerName
    $ws.Range("C$r").Value = $totalPrice
    $ws.Range("D$r").Value = $quantity
    $ws.Range("E$r").Value = $paymentMode
}

# The source data (Total Price / Quantity) are stored as text in this
# report, not numbers -- force Text format on the whole body range before
# writing so Excel doesn't auto-coerce numeric-looking strings, then drop
# back to the Normal style so no stray formatting is left behind.
$bodyRange = $ws.Range("A2:E43")
$bodyRange.NumberFormat = "@"

Set-RowData 2 '63dfe8107771d34b50ec1f77' 'het' '360' '2' 'KOT'
Set-RowData 3 '63e07eaf09901c4fc0bb277e' '' '420' '3' 'CASH'
Set-RowData 4 '63e07e8109901c4fc0bb2751' '' '400' '2' 'CASH'
Set-RowData 5 '63e07f3509901c4fc0bb27fa' '' '280' '2' 'CASH'
Set-RowData 6 '63e0b01b09901c4fc0bb2877' '' '300' '2' 'CASH'
Set-RowData 7 '63e1d5998db3f87bb229b9f2' 'Het B. Patel' '1' '1' 'KOT'
Set-RowData 8 '63e1d7198db3f87bb229bb25' 'Het B. Patel' '400' '2' 'KOT'
Set-RowData 9 '63e1ffe22c5b1158741302ca' '' '360' '3' 'CASH'
Set-RowData 10 '63e203232c5b115874130468' '' '610' '6' 'CASH'
Set-RowData 11 '63e22be88db3f87bb229bb74' 'Het B. Patel' '230' '1' 'KOT'
Set-RowData 12 '63e22d478db3f87bb229bc52' 'Ayushi' '10' '10' 'KOT'
Set-RowData 13 '63e22f0f8db3f87bb229bf35' 'Ayushi' '21' '1' 'KOT'
Set-RowData 14 '63e22f5e8db3f87bb229bf95' 'Ayushi' '60' '1' 'KOT'
Set-RowData 15 '63e22f708db3f87bb229c035' 'Ayushi' '800' '4' 'KOT'
Set-RowData 16 '63e22f938db3f87bb229c08a' 'Ayushi' '2600' '13' 'KOT'
Set-RowData 17 '63e3b550b7feef2bc93c71a7' 'Het B. Patel' '30' '1' 'KOT'
Set-RowData 18 '63e40de811e4eb3328e954f0' 'Het B. Patel' '1630' '8' 'KOT'
Set-RowData 19 '63e4101e11e4eb3328e95562' 'Het B. Patel' '200' '1' 'KOT'
Set-RowData 20 '63e4882611e4eb3328e964db' 'kandarp shah' '230' '1' 'KOT'
Set-RowData 21 '63e486d511e4eb3328e96458' 'kandarp shah' '1080' '5' 'KOT'
Set-RowData 22 '63e48b7311e4eb3328e965ad' 'kandarp shah' '200' '1' 'KOT'
Set-RowData 23 '63e4926111e4eb3328e96993' 'Het B. Patel' '200' '1' 'KOT'
Set-RowData 24 '63e5100c05861c20302bf08d' 'het' '400' '2' 'KOT'
Set-RowData 25 '63e51b274bde7a495d8ec15e' 'Het B. Patel' '60' '1' 'KOT'
Set-RowData 26 '63e53cf2bc76bd300c9dc05c' 'Het B. Patel' '1' '1' 'KOT'
Set-RowData 27 '63e551ac4bde7a495d8ec7d5' 'Het B. Patel' '1' '1' 'KOT'
Set-RowData 28 '63e552041159a16eb41ee78a' '' '280' '2' 'CASH'
Set-RowData 29 '63e55cae4bde7a495d8ec7fe' 'Het B. Patel' '340' '3' 'KOT'
Set-RowData 30 '63e55ce34bde7a495d8ec81a' 'Het B. Patel' '230' '1' 'KOT'
Set-RowData 31 '63e55e954bde7a495d8ec865' 'Het B. Patel' '120' '2' 'KOT'
Set-RowData 32 '63e5627994c22c3ae2c6225a' 'Het B. Patel' '1' '1' 'ONLINE'
Set-RowData 33 '63e563a315af2b3b20982b1a' 'Het B. Patel' '1' '1' 'ONLINE'
Set-RowData 34 '63e565164bde7a495d8ee7ea' 'kandarp shah' '220' '1' 'KOT'
Set-RowData 35 '63e565688dc3223c9284e885' 'Het B. Patel' '1' '1' 'ONLINE'
Set-RowData 36 '63e5d4984bde7a495d8ee8a4' 'Armin' '261' '3' 'KOT'
Set-RowData 37 '63e61a35a743527e2a29a20d' 'Ayushi' '21' '1' 'KOT'
Set-RowData 38 '63e650aca743527e2a29a3b2' 'Armin' '60' '1' 'KOT'
Set-RowData 39 '63e6504ca743527e2a29a373' 'Ayushi' '140' '1' 'KOT'
Set-RowData 40 '63e654a6fbb9f7b872e24a14' '' '27' '2' 'CASH'
Set-RowData 41 '63e6837804b4c646f8cff8ef' 'Het B. Patel' '300' '30' 'KOT'
Set-RowData 42 '63e691d11ecee3b620a9c9be' 'Het B. Patel' '3600' '30' 'KOT'
Set-RowData 43 '63e6846155206a37cc50b08b' 'Het B. Patel' '3600' '30' 'KOT'

$bodyRange.Style = "Normal"
